# Updated cryptos list on Fri Sep 15 11:45:38 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column holds numeric-looking text (e.g. "213.29", "4.10").
# Force it to Text format before writing so the Excel value setter keeps
# the literal string instead of reinterpreting/renormalizing it as a Number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "26.613.13"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "1.630.06"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "213.26"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("E6").Value = "  +2.78%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").Value = "19.19"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").Value = "1.857.65"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.641.41"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "63.54"
$ws.Range("D17").Value = "26.589.31"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").Value = "215.54"
$ws.Range("E19").Value = "  +6.19%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").Value = "9.34"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +5.35%  "
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("E28").Value = "  +3.82%  "
$ws.Range("D29").Value = "15.50"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("D30").Value = "0.0503"
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("D32").Value = "3.29"
$ws.Range("E32").Value = "  +2.94%  "
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "1.222.81"
$ws.Range("E36").Value = "  +5.59%  "
$ws.Range("E37").Value = "  +5.36%  "
$ws.Range("D38").Value = "0.802"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "0.499"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("D42").Value = "0.796"
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("D43").Value = "5.33"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").Value = "1.765.19"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "92.88"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "55.10"
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "0.0513"
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("D50").Value = "7.60"
$ws.Range("E50").Value = "  +3.41%  "
$ws.Range("E51").Value = "  +0.09%  "
